$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the new timesheet entry on row 7
$ws.Range("A7").Value = "Wk [12] Monday 28.5.18"
$ws.Range("B7").Value = "1700 - 2100"
$ws.Range("C7").Value = 4
$ws.Range("D7").Value = "Implementing Priority Queues"

# Move the active selection to C8, matching the saved sheet view state
$ws.Range("C8").Select()

$wb.Save()
